$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 80.5220498566652
$ws.Range("B2").Value = 72.7148740564601
$ws.Range("C2").Value = 88.3292256568703
$ws.Range("D2").Value = 68.5820048356347
$ws.Range("E2").Value = 92.4620948776958

$ws.Range("A3").Value = 80.5220498566652
$ws.Range("B3").Value = 70.5239865382267
$ws.Range("C3").Value = 90.5201131751037
$ws.Range("D3").Value = 65.2313315226252
$ws.Range("E3").Value = 95.8127681907052

$ws.Range("A4").Value = 80.5220498566652
$ws.Range("B4").Value = 68.733472916945
$ws.Range("C4").Value = 92.3106267963854
$ws.Range("D4").Value = 62.4929772449733
$ws.Range("E4").Value = 98.5511224683571

$ws.Range("A5").Value = 80.5220498566652
$ws.Range("B5").Value = 67.1811420012042
$ws.Range("C5").Value = 93.8629577121262
$ws.Range("D5").Value = 60.1188919808762
$ws.Range("E5").Value = 100.925207732454

$ws.Range("A6").Value = 80.5220498566652
$ws.Range("B6").Value = 65.7915000179084
$ws.Range("C6").Value = 95.252599695422
$ws.Range("D6").Value = 57.9936179678035
$ws.Range("E6").Value = 103.050481745527

$ws.Range("A7").Value = 80.5220498566652
$ws.Range("B7").Value = 64.5221006474866
$ws.Range("C7").Value = 96.5219990658438
$ws.Range("D7").Value = 56.052239161897
$ws.Range("E7").Value = 104.991860551433

$ws.Range("A8").Value = 80.5220498566652
$ws.Range("B8").Value = 63.346263096214
$ws.Range("C8").Value = 97.6978366171165
$ws.Range("D8").Value = 54.253950810398
$ws.Range("E8").Value = 106.790148902932

$ws.Range("A9").Value = 80.5220498566652
$ws.Range("B9").Value = 62.2459198783386
$ws.Range("C9").Value = 98.7981798349919
$ws.Range("D9").Value = 52.5711210783363
$ws.Range("E9").Value = 108.472978634994

$ws.Range("A10").Value = 80.5220498566652
$ws.Range("B10").Value = 61.2081635794581
$ws.Range("C10").Value = 99.8359361338723
$ws.Range("D10").Value = 50.9840097788617
$ws.Range("E10").Value = 110.060089934469

$ws.Range("A11").Value = 80.5220498566652
$ws.Range("B11").Value = 60.2233927761319
$ws.Range("C11").Value = 100.820706937199
$ws.Range("D11").Value = 49.4779328019858
$ws.Range("E11").Value = 111.566166911345

$ws.Range("A12").Value = 80.5220498566652
$ws.Range("B12").Value = 59.2842355760646
$ws.Range("C12").Value = 101.759864137266
$ws.Range("D12").Value = 48.0416158113611
$ws.Range("E12").Value = 113.002483901969

$ws.Range("A13").Value = 80.5220498566652
$ws.Range("B13").Value = 58.3848858205589
$ws.Range("C13").Value = 102.659213892772
$ws.Range("D13").Value = 46.6661790535378
$ws.Range("E13").Value = 114.377920659793

Write-Output "Updated range A2:E13 with new forecast values."
